# Add a new "2023" data column (T) to the malaria-incidence table, mirroring
# the existing "2022" column (S) for formatting, then set the 2023 figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: year label for the new column
$ws.Range("S3").Copy($ws.Range("T3"))
$ws.Range("T3").Value = 2023

# Row 4 (Кыргызская Республика / Kyrgyz Republic total)
$ws.Range("S4").Copy($ws.Range("T4"))
$ws.Range("T4").Value = 0

# Rows 5-11 (oblasts) have no 2023 figures yet -> "-" marker,
# matching the style already used for the "-" cells in column S
$ws.Range("S5").Copy($ws.Range("T5"))
$ws.Range("S5").Copy($ws.Range("T6"))
$ws.Range("S5").Copy($ws.Range("T7"))
$ws.Range("S5").Copy($ws.Range("T8"))
$ws.Range("S5").Copy($ws.Range("T9"))
$ws.Range("S5").Copy($ws.Range("T10"))
$ws.Range("S5").Copy($ws.Range("T11"))

# Row 12 (Chui oblast) has an actual 2023 figure
$ws.Range("S12").Copy($ws.Range("T12"))
$ws.Range("T12").Value = 0.001731197036190674

# Row 13 (Bishkek / Osh city row) -> "-" marker
$ws.Range("S13").Copy($ws.Range("T13"))

# Columns A:C were narrowed slightly as part of this edit
$ws.Columns("A:C").ColumnWidth = 32.59
